# Updates '想去人数' (F column) counters and one '最低票价' (G column)
# value (from the non-sellable placeholder text to a numeric price) across
# the four sheets of the workbook, matching the refreshed scrape output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1569
$ws.Range("F5").Value = 8847
$ws.Range("F7").Value = 99
$ws.Range("F8").Value = 1245
$ws.Range("F10").Value = 444
$ws.Range("F11").Value = 593
$ws.Range("F13").Value = 119
$ws.Range("F14").Value = 283
$ws.Range("F16").Value = 50
$ws.Range("F17").Value = 1425
$ws.Range("F19").Value = 566
$ws.Range("F21").Value = 1335
$ws.Range("F23").Value = 211
$ws.Range("F25").Value = 73
$ws.Range("F26").Value = 21
$ws.Range("F28").Value = 284
$ws.Range("F30").Value = 7
$ws.Range("F32").Value = 212
$ws.Range("F33").Value = 180
$ws.Range("F36").Value = 598
$ws.Range("F38").Value = 120
$ws.Range("F41").Value = 468
$ws.Range("F43").Value = 669

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 63
$ws.Range("G23").Value = 180
$ws.Range("F24").Value = 921
$ws.Range("F25").Value = 12
$ws.Range("F27").Value = 193
$ws.Range("F33").Value = 16

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 873
$ws.Range("F5").Value = 736
$ws.Range("F6").Value = 269
$ws.Range("F7").Value = 134
$ws.Range("F8").Value = 1985
$ws.Range("F9").Value = 2997

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1569
$ws.Range("F3").Value = 873
$ws.Range("F4").Value = 736
$ws.Range("F6").Value = 8848
$ws.Range("F7").Value = 269
$ws.Range("F8").Value = 134
$ws.Range("F10").Value = 2997
$ws.Range("F12").Value = 99
$ws.Range("F13").Value = 1245
$ws.Range("F16").Value = 593
$ws.Range("F17").Value = 283
$ws.Range("F18").Value = 50
$ws.Range("F19").Value = 1425
$ws.Range("F21").Value = 566
$ws.Range("F23").Value = 1335
$ws.Range("F24").Value = 211
$ws.Range("F26").Value = 284
$ws.Range("F27").Value = 284
$ws.Range("F30").Value = 63
$ws.Range("F31").Value = 921
$ws.Range("F32").Value = 212
$ws.Range("F33").Value = 12
$ws.Range("F35").Value = 598
$ws.Range("F37").Value = 120
$ws.Range("F41").Value = 468
$ws.Range("F42").Value = 669
